$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new terminal (100000) at the top of the list: shift the
# existing data rows (2..51) down to (3..52), values and formats together,
# working bottom-up so nothing gets overwritten before it is copied.
for ($r = 51; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":G" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":G" + ($r + 1))
    $src.Copy($dst)
}

# Copy the number/text formatting from the row below (now row 3, the old row 2)
# onto the new row 2 for columns A:E, but leave F:G with default (no) style,
# matching how this particular record was entered without the "00000" format.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F2:G2").Style = "Normal"

# Populate the new row's data
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = 354688089433324
$ws.Range("C2").Value = 117212608635
$ws.Range("D2").Value = 8938111000002499584
$ws.Range("E2").Value = 220117701200359
$ws.Range("F2").Value = 100000
$ws.Range("G2").Value = 8657

# Renumber the R.B. (sequence number) column for every data row as text,
# since the new record bumped every existing row's ordinal by one.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = [string]($r - 1)
}

# Add the extra ICCID note in column K next to the 7th original record (now row 8)
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "8938111000002503604"

# Move the active selection to match the saved workbook state
$ws.Range("D11").Select()
